# Generate Report for Handback
# Updates timestamp/status cells that get refreshed when the handback
# report is regenerated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-22 11:39:28"
$wsOverview.Range("G5").Value = "2016-08-22 11:39:28"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-22 11:39:23"
$wsZhCn.Range("H5").Value = "2016-08-22 11:39:23"
$wsZhCn.Range("K3").Value = "2016-08-22 11:39:41"
$wsZhCn.Range("K5").Value = "2016-08-22 11:39:41"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-22 11:39:28"
$wsDeDe.Range("H5").Value = "2016-08-22 11:39:28"
$wsDeDe.Range("K3").Value = "2016-08-22 11:39:48"
$wsDeDe.Range("K5").Value = "2016-08-22 11:39:48"
